$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-10-13 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-14 Monday", 2) | Out-Null
$d.Content.Find.Execute("74×61=", $true, $false, $false, $false, $false, $true, 1, $false, "68×84=", 2) | Out-Null
$d.Content.Find.Execute("86×94=", $true, $false, $false, $false, $false, $true, 1, $false, "70×65=", 2) | Out-Null
$d.Content.Find.Execute("35×92=", $true, $false, $false, $false, $false, $true, 1, $false, "46×13=", 2) | Out-Null
$d.Content.Find.Execute("98×58=", $true, $false, $false, $false, $false, $true, 1, $false, "13×79=", 2) | Out-Null
$d.Content.Find.Execute("65×54=", $true, $false, $false, $false, $false, $true, 1, $false, "23×93=", 2) | Out-Null
$d.Content.Find.Execute("58×81=", $true, $false, $false, $false, $false, $true, 1, $false, "69×83=", 2) | Out-Null
$d.Content.Find.Execute("53×96=", $true, $false, $false, $false, $false, $true, 1, $false, "60×46=", 2) | Out-Null
$d.Content.Find.Execute("88×56=", $true, $false, $false, $false, $false, $true, 1, $false, "74×32=", 2) | Out-Null
$d.Content.Find.Execute("39×42=", $true, $false, $false, $false, $false, $true, 1, $false, "59×29=", 2) | Out-Null
$d.Content.Find.Execute("86×79=", $true, $false, $false, $false, $false, $true, 1, $false, "20×94=", 2) | Out-Null
$d.Content.Find.Execute("13×95=", $true, $false, $false, $false, $false, $true, 1, $false, "97×53=", 2) | Out-Null
$d.Content.Find.Execute("19×94=", $true, $false, $false, $false, $false, $true, 1, $false, "92×67=", 2) | Out-Null
$d.Content.Find.Execute("82×40=", $true, $false, $false, $false, $false, $true, 1, $false, "48×95=", 2) | Out-Null
$d.Content.Find.Execute("48×53=", $true, $false, $false, $false, $false, $true, 1, $false, "12×22=", 2) | Out-Null
$d.Content.Find.Execute("30×25=", $true, $false, $false, $false, $false, $true, 1, $false, "12×26=", 2) | Out-Null
$d.Content.Find.Execute("34×94=", $true, $false, $false, $false, $false, $true, 1, $false, "35×89=", 2) | Out-Null
$d.Content.Find.Execute("84×78=", $true, $false, $false, $false, $false, $true, 1, $false, "37×79=", 2) | Out-Null
$d.Content.Find.Execute("75×60=", $true, $false, $false, $false, $false, $true, 1, $false, "68×17=", 2) | Out-Null
$d.Content.Find.Execute("33×68=", $true, $false, $false, $false, $false, $true, 1, $false, "51×32=", 2) | Out-Null
$d.Content.Find.Execute("75×63=", $true, $false, $false, $false, $false, $true, 1, $false, "30×53=", 2) | Out-Null
$d.Content.Find.Execute("84×69=", $true, $false, $false, $false, $false, $true, 1, $false, "95×86=", 2) | Out-Null
$d.Content.Find.Execute("55×62=", $true, $false, $false, $false, $false, $true, 1, $false, "33×74=", 2) | Out-Null
$d.Content.Find.Execute("53×23=", $true, $false, $false, $false, $false, $true, 1, $false, "31×96=", 2) | Out-Null
$d.Content.Find.Execute("16×14=", $true, $false, $false, $false, $false, $true, 1, $false, "83×78=", 2) | Out-Null
$d.Content.Find.Execute("25×60=", $true, $false, $false, $false, $false, $true, 1, $false, "68×84=", 2) | Out-Null
